$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.126.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.279.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '155.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15,427.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '94.90'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.17%  '

$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.492'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '35.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.91%  '

$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("E13").Value = '  -1.92%  '

$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.632.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.18%  '

$ws.Range("E16").Value = '  +1.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.282.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.796'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.043.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0916'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '243.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.48%  '

$ws.Range("E30").Value = '  +1.10%  '

$ws.Range("E31").Value = '  +1.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.48%  '

$ws.Range("E33").Value = '  +3.23%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0752'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("E36").Value = '  +0.99%  '

$ws.Range("E37").Value = '  +3.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.67%  '

$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.014.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("E45").Value = '  +11.47%  '

$ws.Range("E46").Value = '  +1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.27%  '

$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.73%  '
